$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.280.26"
$ws.Range("E2").Value = "  +0.28%  "
$ws.Range("D3").Value = "1.680.76"
$ws.Range("E3").Value = "  +0.61%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "'218.56"
$ws.Range("E5").Value = "  +0.47%  "
$ws.Range("D6").Value = "'0.5264"
$ws.Range("E6").Value = "  +2.81%  "
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("E8").Value = "  +1.90%  "
$ws.Range("D9").Value = "'0.06438"
$ws.Range("E9").Value = "  +0.90%  "
$ws.Range("E10").Value = "  +2.64%  "
$ws.Range("E11").Value = "  +1.69%  "
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").Value = "'4.548"
$ws.Range("E12").Value = "  +0.17%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.678.79"
$ws.Range("E13").Value = "  +0.42%  "
$ws.Range("D14").Value = "'0.5814"
$ws.Range("E14").Value = "  -0.04%  "
$ws.Range("D15").Value = "'0.000008505"
$ws.Range("E15").Value = "  -1.60%  "
$ws.Range("D16").Value = "'64.53"
$ws.Range("E16").Value = "  +0.17%  "
$ws.Range("D17").Value = "26.334.53"
$ws.Range("E17").Value = "  +0.49%  "
$ws.Range("D18").Value = "'4.937"
$ws.Range("E18").Value = "  +0.01%  "
$ws.Range("E19").Value = "  +0.07%  "
$ws.Range("E20").Value = "  +0.14%  "
$ws.Range("D21").Value = "'189.76"
$ws.Range("E21").Value = "  +0.17%  "
$ws.Range("D22").Value = "'6.219"
$ws.Range("E22").Value = "  +0.12%  "
$ws.Range("D23").Value = "'1.008"
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("D24").Value = "'145.32"
$ws.Range("D25").Value = "'7.768"
$ws.Range("E25").Value = "  +1.71%  "
$ws.Range("D26").Value = "'0.1248"
$ws.Range("E26").Value = "  +5.98%  "
$ws.Range("E27").Value = "  +1.12%  "
$ws.Range("D28").Value = "'0.06630"
$ws.Range("E28").Value = "  +10.98%  "
$ws.Range("E29").Value = "  +5.60%  "
$ws.Range("E30").Value = "  +0.47%  "
$ws.Range("D31").Value = "'3.592"
$ws.Range("E31").Value = "  +1.96%  "
$ws.Range("D32").Value = "'3.578"
$ws.Range("E32").Value = "  +1.51%  "
$ws.Range("D33").Value = "'1.662"
$ws.Range("E33").Value = "  +1.13%  "
$ws.Range("D34").Value = "'1.029"
$ws.Range("E34").Value = "  +1.45%  "
$ws.Range("D35").Value = "'0.6224"
$ws.Range("E35").Value = "  +3.06%  "
$ws.Range("E36").Value = "  +0.95%  "
$ws.Range("E37").Value = "  +2.60%  "
$ws.Range("D38").Value = "'6.443"
$ws.Range("E38").Value = "  +5.82%  "
$ws.Range("D39").Value = "1.108.90"
$ws.Range("E39").Value = "  +2.95%  "
$ws.Range("E40").Value = "  +0.14%  "
$ws.Range("D41").Value = "'0.8797"
$ws.Range("E41").Value = "  +1.28%  "
$ws.Range("E42").Value = "  +0.47%  "
$ws.Range("D43").Value = "'100.72"
$ws.Range("E43").Value = "  +0.37%  "
$ws.Range("D44").Value = "1.838.08"
$ws.Range("E44").Value = "  +0.96%  "
$ws.Range("E45").Value = "  +1.67%  "
$ws.Range("D46").Value = "'57.00"
$ws.Range("E46").Value = "  +1.26%  "
$ws.Range("E47").Value = "  -0.28%  "
$ws.Range("D48").Value = "'8.159"
$ws.Range("E48").Value = "  +1.04%  "
$ws.Range("D49").Value = "'0.05278"
$ws.Range("E49").Value = "  +1.30%  "
$ws.Range("D50").Value = "'0.4300"
$ws.Range("E51").Value = "  +3.35%  "
